# Updated cryptos list refresh: coin prices, 1h volume %, and two list re-ranks
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (coin name, link, %-volume, and non-numeric-looking prices)
# can be set directly; Excel keeps them as text because they are not valid numbers.
$ws.Range("D2").Value = "62.928.83"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "2.682.54"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("E9").Value = "  -3.25%  "
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.368"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "3.156.53"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").Value = "62.821.37"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "2.683.30"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.73%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.507"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("D28").Value = "0.0₃0859"
$ws.Range("E28").Value = "  -6.47%  "
$ws.Range("E29").Value = "  +3.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "350.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.963"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.16%  "
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.87%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.617"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0561"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0973"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("E50").Value = "  -3.44%  "
$ws.Range("D51").Value = "2.097.63"
$ws.Range("E51").Value = "  -1.77%  "
